$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column Q (values + formatting/styles) into the new column R so
# the new column inherits the exact same per-row cell styles as column Q.
$ws.Columns("Q:Q").Copy() | Out-Null
$ws.Columns("R:R").Insert(-4161) | Out-Null

# Overwrite the copied values in column R with the 2021 figures while
# keeping the formatting that was just duplicated from column Q.
$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 5.8
$ws.Range("R6").Value = 4.7
$ws.Range("R7").Value = 1.6
$ws.Range("R8").Value = 12.9
$ws.Range("R9").Value = 10.199999999999999
$ws.Range("R10").Value = 4.2
$ws.Range("R11").Value = 3.3
$ws.Range("R12").Value = 15.2
$ws.Range("R13").Value = 2.4
$ws.Range("R14").Value = 0.6

# Update the active selection to match the author's final cursor position.
$ws.Range("T9").Select() | Out-Null
